$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ad deletion implemented - award full score (5) for the related
# self-evaluation rows: "Edit Inactive Ads" (C26) and "Delete Ad" (C28).
$ws.Range("C26").Value = 5
$ws.Range("C28").Value = 5

# Reflect the reviewer's current scroll/selection position in the sheet view.
$ws.Range("G26").Select()
